$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F (reuse the bold header formatting from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamp values for rows 2-41 (as text, matching inlineStr in source)
$timestamps = @(
    "2021-10-05 13:42:14.287295",
    "2021-10-05 13:42:14.287306",
    "2021-10-05 13:42:14.287309",
    "2021-10-05 13:42:14.287312",
    "2021-10-05 13:42:14.287316",
    "2021-10-05 13:42:14.287318",
    "2021-10-05 13:42:14.287321",
    "2021-10-05 13:42:14.287323",
    "2021-10-05 13:42:14.287326",
    "2021-10-05 13:42:14.287329",
    "2021-10-05 13:42:14.287332",
    "2021-10-05 13:42:14.287334",
    "2021-10-05 13:42:14.287337",
    "2021-10-05 13:42:14.287340",
    "2021-10-05 13:42:14.287343",
    "2021-10-05 13:42:14.287345",
    "2021-10-05 13:42:14.287348",
    "2021-10-05 13:42:14.287351",
    "2021-10-05 13:42:14.287353",
    "2021-10-05 13:42:14.287356",
    "2021-10-05 13:42:14.287358",
    "2021-10-05 13:42:14.287361",
    "2021-10-05 13:42:14.287363",
    "2021-10-05 13:42:14.287366",
    "2021-10-05 13:42:14.287369",
    "2021-10-05 13:42:14.287372",
    "2021-10-05 13:42:14.287375",
    "2021-10-05 13:42:14.287377",
    "2021-10-05 13:42:14.287380",
    "2021-10-05 13:42:14.287382",
    "2021-10-05 13:42:14.287385",
    "2021-10-05 13:42:14.287388",
    "2021-10-05 13:42:14.287391",
    "2021-10-05 13:42:14.287394",
    "2021-10-05 13:42:14.287396",
    "2021-10-05 13:42:14.287398",
    "2021-10-05 13:42:14.287401",
    "2021-10-05 13:42:14.287404",
    "2021-10-05 13:42:14.287407",
    "2021-10-05 13:42:14.287409"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
